# Add a new worksheet "工作表3" with a small A/B/C/D score table, update the
# selection/active-tab bookkeeping left behind by Excel, matching the
# "Add files via upload" commit.

$wb = $excel.ActiveWorkbook

# --- 1. Touch 工作表1: it was the previously-selected tab (selection C3);
#        the new save leaves it on cell F9 and no longer the active tab. ---
$ws1 = $wb.Worksheets.Item("工作表1")
$ws1.Activate()
$ws1.Range("F9").Select()

# --- 2. Add the new worksheet 工作表3 after the existing 工作表2 ---
$ws2 = $wb.Worksheets.Item("工作表2")
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "工作表3"

# --- 3. Fill in the data (matches sharedStrings A/B/C/D + scores) ---
$ws3.Range("B2").Value = "A"
$ws3.Range("C2").Value = 20
$ws3.Range("B3").Value = "B"
$ws3.Range("C3").Value = 40
$ws3.Range("B4").Value = "C"
$ws3.Range("C4").Value = 5
$ws3.Range("B5").Value = "D"
$ws3.Range("C5").Value = 35

# --- 4. Match the page setup used by the other sheets ---
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

# --- 5. 工作表3 becomes the active sheet/tab, selection resting on E6 ---
$ws3.Range("E6").Select()
